$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 2600
    $ws.Range("F4").Value = 523
    $ws.Range("F6").Value = 6581
    $ws.Range("F7").Value = 424
    $ws.Range("F8").Value = 8
    $ws.Range("F9").Value = 7
    $ws.Range("F10").Value = 10
    $ws.Range("F11").Value = 3

    $ws.Range("C12").Value = "合肥·梦时空SPO1动漫展（取消）"
    $ws.Range("G12").Value = "不可售"
}
